$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------
# Sheet "Judge": add a new "undergrad entry year" column and replace
# the sample judge row with the new judge's data.
# ----------------------------------------------------------------
$judge = $wb.Worksheets.Item("Judge")

# Insert a new column before the old "QQ" column (old column B), which
# pushes QQ / phone / resume one column to the right.
$judge.Columns.Item(2).Insert()

$judge.Range("B2").Value = "评委本科入学年份"
$judge.Range("E2").Value = "评委履历"
# the resume header used to be bold (like the rest of row 2); after the
# insert it keeps that formatting, but the final sheet shows it un-bolded.
$judge.Range("E2").Font.Bold = $false

$judge.Range("A3").Value = "高子柳"
$judge.Range("B3").Value = "2015"
$judge.Range("C3").Value = "1010678911"
$judge.Range("D3").Value = "13880551583"
$judge.Range("E3").Value = "很厉害，超厉害"

# Approximate the auto-fit column widths for the new/resized columns.
$judge.Columns.Item(2).ColumnWidth = 19.02
$judge.Columns.Item(5).ColumnWidth = 16.6

$judge.Activate()
$judge.Range("C9").Select()

# ----------------------------------------------------------------
# Sheet "Basic": update the submission-timing row (row 3)
# ----------------------------------------------------------------
$basic = $wb.Worksheets.Item("Basic")

# A3 keeps being a hyperlink, but both its address and display text move
# to the new test-link URL.
foreach ($h in $basic.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$3') {
        $h.Address = "https://www.wjx.cn/jq/42340427.aspx"
    }
}
$basic.Range("A3").Value = "https://www.wjx.cn/jq/42340427.aspx"

$basic.Range("B3").Value = "2019-07-03  00:01:00"
$basic.Range("C3").Value = "2019-07-07  16:53:00"

# D3 ("official link") no longer has a value or a hyperlink.
foreach ($h in $basic.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.Delete()
    }
}
$basic.Range("D3").ClearContents()

$basic.Range("E3").Value = "2019-07-07 16:53:00"
$basic.Range("F3").Value = "2019-07-07 16:55:00"

$basic.Activate()
$basic.Range("D3").Select()
